# Add more exploration of variables: populate tillage (E) and n_mgmt (F)
# columns on the "over-years" sheet, shifting the old notes/notes2 columns
# out to G/H, and fill in the per-study notes/notes2 detail that goes with
# the new n_mgmt column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("over-years")

# --- Header row -----------------------------------------------------
# E1 "tillage" already correct. Insert "n_mgmt" after it; push the old
# "notes"/"notes2" headers one column to the right (F,G -> G,H).
$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("F1").Value2 = "n_mgmt"

# --- Siefert et al. 2017 block (rows 2-6): tillage "varies", n_mgmt
# "farmer choice" ------------------------------------------------------
$ws.Range("E2:E6").Value = "varies"
$ws.Range("F2:F6").Value = "farmer choice"

# --- Gentry et al. 2013 block (rows 7-10) ---------------------------
$ws.Range("E7:E10").Value = "chisel-plow fall, cultivated spring"
$ws.Range("F7:F10").Value = "multiple rates"

# Move the old per-row "notes" (range of x%) from F8:F10 out to G8:G10,
# and retire the old G7 "notes2" value (cleared; row 7 has no G value any
# more), replacing the notes2 column content (H7:H10) with new per-row
# detail.
$ws.Range("G7").ClearContents()
$ws.Range("G8").Value = "range of 4-10%"
$ws.Range("G9").Value = "range of 12-16%"
$ws.Range("G10").Value = "range of 17-22%"

$ws.Range("H7").Value = "year is confounded with #years in corn"
$ws.Range("H8").Value = "Yield gap determined from regression equations at AONR, not actually measured. "
$ws.Range("H9").Value = "Used different N rates each year"
$ws.Range("H10").Value = "7th year corn received smallest N rate"

# --- Crookston et al. 1991 block (rows 11-15) ------------------------
$ws.Range("E11:E15").Value = "moldboard fall, spring disk"
$ws.Range("F11:F15").Value = "aonr"

# --- Porter et al. 1997 blocks (rows 16-33) --------------------------
$ws.Range("E16:E33").Value = "moldboard fall, spring disk"
$ws.Range("F16:F33").Value = "aonr"

# --- Meese et al. 1991 block (rows 34-36) ----------------------------
$ws.Range("E34:E36").Value = "moldboard fall, spring disk"
$ws.Range("F34:F36").Value = "aonr"

# --- Column widths: new n_mgmt column (F) matches D:E width ----------
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# --- Selection matches where the author last clicked ------------------
$ws.Range("G7").Select()
